$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 34: paper added to the ARI reference list (March 2021 - June 2021)
$ws.Range("A34").Value = 'Across-subjects multiple baseline trial of exposure-based cognitive-behavioral therapy for severe irritability: a study protocol'
$ws.Range("B34").Value = 'Reut Naim, Katharina Kircanski, Andrea Gold,  Ramaris E German, Mollie Davis, Samantha Perlstein, Michal Clayton, Olga Revzina, and Melissa A Brotman'
$ws.Range("C34").Value = 'BMJ Open'
$ws.Range("D34").Value = 2021
$ws.Range("E34").Value = 'https://dx.doi.org/10.1136%2Fbmjopen-2020-039169'
$ws.Range("F34").Value = @'
Introduction Irritability is defined as a tendency towards anger in response to frustration. Clinically, impairing irritability is a significant public health problem. There is a need for mechanism-based psychotherapies targeting severe irritability as it manifests in the context of disruptive mood dysregulation disorder (DMDD). This study protocol describes a randomised multiple baseline design testing the preliminary efficacy of a new treatment, exposure-based cognitive-behavioral therapy for severe irritability in youth, which also integrates components of parent management training. We will investigate associations of this intervention with primary clinical measures, as well as ecological momentary assessment measures.
Methods and analysis Forty youth will be enrolled. Participants, aged 8–17 years, must present at least one of two core symptoms of DMDD: abnormal mood or increased reactivity to negative emotional stimuli, with severe impairment in one domain (home, school, peers) and moderate in another, or moderate impairment in at least two domains. Each participant is randomised to a 2-week, 4-week or 6-week baseline observation period, followed by 12 active treatment sessions. Clinical ratings are conducted at baseline, biweekly (clinician), weekly (parent/child) throughout treatment, post-treatment, and 3-month and 6-month follow-up (clinician). Clinician ratings on the Affective Reactivity Index and Clinical Global Impressions-Improvement scale for DMDD are our primary outcome measures. Secondary outcome measures include parent and child reports of irritability. Post hoc additional symptom measures include clinician, parent and self-ratings of depression, anxiety and overall functional impairment. Prospective, digitally based event sampling of symptoms is acquired for a week pre-treatment, mid-treatment and post-treatment. Based on our pathophysiological model of irritability implicating frustrative non-reward, aberrant threat processing and instrumental learning, we probe these three brain-based targets using functional MRI paradigms to assess target engagement.
Ethics and dissemination The research project and all related materials were submitted and approved by the appropriate Institutional Review Board (IRB) of the National Institute of Mental Health (NIMH).
Trial registration numbers NCT02531893 and NCT00025935.
http://creativecommons.org/licenses/by-nc/4.0/
This is an open access article distributed in accordance with the Creative Commons Attribution Non Commercial (CC BY-NC 4.0) license, which permits others to distribute, remix, adapt, build upon this work non-commercially, and license their derivative works on different terms, provided the original work is properly cited, appropriate credit is given, any changes made indicated, and the use is non-commercial. See: http://creativecommons.org/licenses/by-nc/4.0/.
'@
$ws.Range("G34").Value = 'ARI used as part of measures'
$ws.Range("H34").Value = 'Treatment'

# I34 (Type_Secondary) is left blank for this paper, but the cell is still
# present in the row like the rest of the table (e.g. I33).
$ws.Range("I34").Font.Bold = $false

# Setting a multi-line Abstract bumps this row to a custom height; re-run
# AutoFit so row 34 keeps the sheet's normal (non-custom) row height, same
# as every other row in the table.
$ws.Rows.Item(34).AutoFit()
